{"js": "// The document contains a single, one-column table of DaCapo/ZGC benchmark\n// stats. This edit:\n//   1. Rewrites the heap-size summary rows (rows 1-3) from \"100\"/\"0\"/\"17\"\n//      to \"0M\"/\"0M\"/\"0M\".\n//   2. Expands a single tab-delimited \"stats line\" row (which used to hold\n//      10 separate w:t runs joined by w:tab) into 10 separate table rows,\n//      one value per row, inserted right after the heap-size rows.\n//   3. Collapses the two remaining tab-delimited multi-run rows down to a\n//      single value each (\"100\" and \"0\").\n//   4. Fills in the previously-empty trailing row with \"17\".\n//\n// Net effect: the table grows from 36 rows to 46 rows, and every row ends\n// up holding exactly one value in a single run.\n\nconst table = context.document.body.tables.getFirst();\ntable.load(\"rowCount\");\nawait context.sync();\n\n// Make room for the 10 newly-introduced rows (one per tab-separated value\n// that used to live in a single cell).\ntable.addRows(Word.InsertLocation.end, 10);\nawait context.sync();\n\n// Final, complete set of per-row values (single column) once every row\n// holds just its own value.\nconst finalValues = [\n  [\"0M\"],\n  [\"0M\"],\n  [\"0M\"],\n  [\"13\"],\n  [\"0.00003\"],\n  [\"0.00004\"],\n  [\"0.00004\"],\n  [\"0.00000\"],\n  [\"0.00003\"],\n  [\"0.00004\"],\n  [\"0.00004\"],\n  [\"0.00046\"],\n  [\"100.0\"],\n  [\"0\"],\n  [\"0.00000\"],\n  [\"0.00000\"],\n  [\"0.00000\"],\n  [\"0.00000\"],\n  [\"0.00000\"],\n  [\"0.00000\"],\n  [\"0.00000\"],\n  [\"0.00000\"],\n  [\"0.0\"],\n  [\"0\"],\n  [\"0.00000\"],\n  [\"0.00000\"],\n  [\"0.00000\"],\n  [\"0.00000\"],\n  [\"0.00000\"],\n  [\"0.00000\"],\n  [\"0.00000\"],\n  [\"0.00000\"],\n  [\"0.0\"],\n  [\"78\"],\n  [\"0.00000\"],\n  [\"0.30094\"],\n  [\"0.03266\"],\n  [\"0.01566\"],\n  [\"0.12622\"],\n  [\"0.13196\"],\n  [\"0.17117\"],\n  [\"2.54776\"],\n  [\"551463.9\"],\n  [\"100\"],\n  [\"0\"],\n  [\"17\"],\n];\n\ntable.values = finalValues;\nawait context.sync();\n", "ps1": "# The document contains a single, one-column table of DaCapo/ZGC benchmark\n# stats. This edit:\n#   1. Rewrites the heap-size summary rows (rows 1-3) from \"100\"/\"0\"/\"17\"\n#      to \"0M\"/\"0M\"/\"0M\".\n#   2. Expands a single tab-delimited \"stats line\" row (which used to hold\n#      10 separate runs joined by tab characters) into 10 separate table\n#      rows, one value per row, inserted right after the heap-size rows.\n#   3. Collapses the two remaining tab-delimited multi-run rows down to a\n#      single value each (\"100\" and \"0\").\n#   4. Fills in the previously-empty trailing row with \"17\".\n#\n# Net effect: the table grows from 36 rows to 46 rows, and every row ends\n# up holding exactly one value in a single run.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Make room for the 10 newly-introduced rows (one per tab-separated value\n# that used to live in a single cell).\nfor ($i = 0; $i -lt 10; $i++) {\n    $t.Rows.Add() | Out-Null\n}\n\n# Final, complete set of per-row values (single column) once every row\n# holds just its own value.\n$finalValues = @(\n    \"0M\",\n    \"0M\",\n    \"0M\",\n    \"13\",\n    \"0.00003\",\n    \"0.00004\",\n    \"0.00004\",\n    \"0.00000\",\n    \"0.00003\",\n    \"0.00004\",\n    \"0.00004\",\n    \"0.00046\",\n    \"100.0\",\n    \"0\",\n    \"0.00000\",\n    \"0.00000\",\n    \"0.00000\",\n    \"0.00000\",\n    \"0.00000\",\n    \"0.00000\",\n    \"0.00000\",\n    \"0.00000\",\n    \"0.0\",\n    \"0\",\n    \"0.00000\",\n    \"0.00000\",\n    \"0.00000\",\n    \"0.00000\",\n    \"0.00000\",\n    \"0.00000\",\n    \"0.00000\",\n    \"0.00000\",\n    \"0.0\",\n    \"78\",\n    \"0.00000\",\n    \"0.30094\",\n    \"0.03266\",\n    \"0.01566\",\n    \"0.12622\",\n    \"0.13196\",\n    \"0.17117\",\n    \"2.54776\",\n    \"551463.9\",\n    \"100\",\n    \"0\",\n    \"17\"\n)\n\nfor ($i = 0; $i -lt $finalValues.Count; $i++) {\n    $cell = $t.Cell($i + 1, 1)\n    $cell.Range.Text = $finalValues[$i]\n}\n"}
